# Correcion a Diebold Mariano y revision de Cap1
# Updates DM_Stat (col C) and P_Value (col D) for rows 2-11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  C = -0.9614084219133069;  D = 0.3431347699486653 },
    @{ Row = 3;  C = -1.030093843559184;   D = 0.3102384205220745 },
    @{ Row = 4;  C = -0.9242416843674782;  D = 0.3618739130672126 },
    @{ Row = 5;  C = -0.7777706181314397;  D = 0.4420830269641263 },
    @{ Row = 6;  C = -0.6541119378806109;  D = 0.5174399601905186 },
    @{ Row = 7;  C = -0.4940317269504514;  D = 0.6244605306126427 },
    @{ Row = 8;  C = -0.5294001909921355;  D = 0.5999677572530495 },
    @{ Row = 9;  C = 0.2963530163120485;   D = 0.7687632408037954 },
    @{ Row = 10; C = 0.04209592419873306;  D = 0.9666684732129343 },
    @{ Row = 11; C = -0.2686676325921321;  D = 0.7898090154478539 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
